$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Data Lançamento" (column D) values from plain years to real dates,
# and format that column as a date.
$ws.Range("D2").Value = [DateTime]"2020-10-20"
$ws.Range("D3").Value = [DateTime]"2021-10-10"
$ws.Range("D4").Value = [DateTime]"2015-05-05"
$ws.Range("D5").Value = [DateTime]"2021-01-12"
$ws.Range("D6").Value = [DateTime]"1995-05-16"
$ws.Range("D2:D6").NumberFormat = "m/d/yyyy"

# Paint entire used columns white first (background clean-up)
$ws.Range("A:L").Interior.ThemeColor = 2

# Green fill for Nome, Descrição, Duração, Data Lançamento, Lançamento, Avaliação
$ws.Range("A1:E6").Interior.Color = 5296274
$ws.Range("H1:H6").Interior.Color = 5296274

# Yellow fill for Categoria, Classificação, Atores, Diretores, Escritores, Tipo de filme
$ws.Range("F1:G6").Interior.Color = 65535
$ws.Range("I1:L6").Interior.Color = 65535

# Move the active selection like in the diff
$ws.Range("B3").Select()
